$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column A (GENE id duplicate) is removed entirely; columns B:F shift left to A:E.
$ws.Range("A1").EntireColumn.Delete()
